$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 360 - this shifts the existing rows 360..390 down to 361..391
$ws.Rows.Item(360).Insert()

# Populate the newly inserted row 360 with the new record
$ws.Range("A360").Value = 5
$ws.Range("B360").Value = "Macroferia Regional de Talca"
$ws.Range("C360").Value = "Maule"
$ws.Range("D360").Value = 45265
$ws.Range("E360").Value = 7
$ws.Range("F360").Value = 100112021
$ws.Range("G360").Value = "Ají"
$ws.Range("H360").Value = "Inferno"
$ws.Range("I360").Value = "Primera"
$ws.Range("J360").Value = 100
$ws.Range("K360").Value = 22000
$ws.Range("L360").Value = 22000
$ws.Range("M360").Value = 22000
$ws.Range("N360").Value = "$/caja 14 kilos"
$ws.Range("O360").Value = "Región del Maule"
$ws.Range("P360").Value = 1571
$ws.Range("Q360").Value = 14
$ws.Range("R360").Value = "Hortaliza"
